# Weekly roll-forward of the "Hortaliza, Vega Monumental Concepción - Perejil"
# logica_diaria subset: every existing observation (rows 32-129, a run of
# paired Primera/Segunda records sorted most-recent-first starting at row 32)
# shifts down by one pair (two rows) to make room for a brand new, most
# recent pair inserted at rows 32-33. The new pair reuses the same
# market/price/volume/origin values that used to sit in rows 32-33, just
# stamped with a newer date (serial 44607 = 2022-02-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")

$firstRow = 32
$lastRow  = 129
$shift    = 2

# Work from the bottom up so the source row for each copy hasn't been
# overwritten yet.
for ($r = $lastRow + $shift; $r -ge $firstRow + $shift; $r--) {
    $src = $r - $shift
    foreach ($c in $cols) {
        $ws.Range($c + $r).Value = $ws.Range($c + $src).Value2
    }
    # Keep the date column's date number-format on the (possibly brand new)
    # destination row in sync with the rest of the column.
    $ws.Range("D" + $r).NumberFormat = $ws.Range("D" + $src).NumberFormat
}

# The two rows at the very top of this block (32-33) keep all of their old
# values/formats; only the date changes, to the new, most recent date.
$ws.Range("D32").Value = 44607
$ws.Range("D33").Value = 44607
